# Feat: Add ${username.id} to replace username with correct id
#
# Replace the hard-coded "/redfish/v1/AccountService/Accounts/testuser"
# endpoint references (used by the PATCH/DELETE rows) with a templated
# endpoint that uses the created user's id: "/redfish/v1/AccountService/Accounts/${testuser.id}"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RedfishCommands")

$oldValue = "/redfish/v1/AccountService/Accounts/testuser"
$newValue = '/redfish/v1/AccountService/Accounts/${testuser.id}'

$usedRange = $ws.UsedRange
for ($r = 1; $r -le $usedRange.Rows.Count; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    if ($cell.Value() -eq $oldValue) {
        $cell.Value = $newValue
    }
}

# Move the active selection as recorded by the saved workbook
$ws.Range("B15").Select()
